# Updates the cryptocurrency price/volume snapshot in columns D (Price)
# and E (Volume(1h)) for the data rows of the active sheet, mirroring a
# refreshed pull of the source feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row => @{ column letter = new text value }
# Only D (Price) and/or E (Volume(1h)) are present for rows whose figures moved.
$updates = @{
    2 = @{ "D" = "57.097.21"; "E" = "  +0.65%  " }
    3 = @{ "D" = "2.422.27"; "E" = "  -1.83%  " }
    4 = @{ "E" = "  +0.07%  " }
    5 = @{ "D" = "488.42"; "E" = "  -0.11%  " }
    6 = @{ "D" = "153.97"; "E" = "  +1.98%  " }
    7 = @{ "E" = "  +19.78%  " }
    8 = @{ "D" = "0.996"; "E" = "  -0.15%  " }
    9 = @{ "D" = "2.443.63"; "E" = "  -1.37%  " }
    10 = @{ "D" = "6.20"; "E" = "  +8.52%  " }
    11 = @{ "D" = "0.1000"; "E" = "  +0.53%  " }
    12 = @{ "E" = "  -0.31%  " }
    13 = @{ "E" = "  +1.32%  " }
    14 = @{ "D" = "2.848.26"; "E" = "  -1.88%  " }
    15 = @{ "D" = "57.146.92"; "E" = "  +0.11%  " }
    16 = @{ "D" = "20.56"; "E" = "  -1.97%  " }
    18 = @{ "D" = "2.439.17"; "E" = "  -1.67%  " }
    19 = @{ "E" = "  +1.60%  " }
    20 = @{ "D" = "324.38"; "E" = "  +1.41%  " }
    21 = @{ "D" = "9.98"; "E" = "  -1.67%  " }
    23 = @{ "D" = "5.91"; "E" = "  +1.17%  " }
    24 = @{ "D" = "57.78"; "E" = "  -0.34%  " }
    25 = @{ "D" = "0.400"; "E" = "  -1.18%  " }
    26 = @{ "D" = "0.995"; "E" = "  -0.45%  " }
    27 = @{ "E" = "  -1.41%  " }
    28 = @{ "D" = "2.533.94"; "E" = "  -2.32%  " }
    29 = @{ "D" = "7.26"; "E" = "  -4.02%  " }
    30 = @{ "D" = "0.0₃0786"; "E" = "  -2.22%  " }
    31 = @{ "E" = "  -0.05%  " }
    32 = @{ "D" = "150.63"; "E" = "  -0.21%  " }
    33 = @{ "D" = "18.68"; "E" = "  +2.41%  " }
    34 = @{ "E" = "  +0.46%  " }
    35 = @{ "E" = "  +1.78%  " }
    36 = @{ "D" = "3.78"; "E" = "  +0.70%  " }
    37 = @{ "E" = "  -0.74%  " }
    38 = @{ "D" = "0.820"; "E" = "  -7.42%  " }
    39 = @{ "D" = "0.101"; "E" = "  +7.47%  " }
    40 = @{ "D" = "285.03"; "E" = "  +8.88%  " }
    41 = @{ "D" = "34.07"; "E" = "  -0.10%  " }
    42 = @{ "E" = "  -0.70%  " }
    43 = @{ "E" = "  +0.62%  " }
    44 = @{ "D" = "0.995"; "E" = "  -0.19%  " }
    45 = @{ "D" = "0.601"; "E" = "  -0.87%  " }
    46 = @{ "E" = "  -4.53%  " }
    47 = @{ "E" = "  -0.07%  " }
    48 = @{ "D" = "0.0227"; "E" = "  -0.43%  " }
    49 = @{ "E" = "  -4.37%  " }
    50 = @{ "D" = "1.895.77"; "E" = "  +2.28%  " }
    51 = @{ "D" = "17.59"; "E" = "  -0.70%  " }
}

$colIndex = @{ "D" = 4; "E" = 5 }

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $cell = $ws.Cells.Item($row, $colIndex[$col])
        # Many of the values look numeric ("488.42", "6.20", "0.1000", ...)
        # but the source column stores them as plain text (leading/trailing
        # zeros, thousand-dot grouping, padded percentages, etc. must survive
        # verbatim). Force text interpretation while writing, then drop back
        # to the default style so no extra formatting is left behind.
        $cell.NumberFormat = "@"
        $cell.Value = $updates[$row][$col]
        $cell.Style = "Normal"
    }
}
